# Update the PSSM score matrix (B2:K21) with the recomputed values from the
# supplemental-figure re-run. Row/column headers (row 1, column A) and the
# cell at A1 are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 20,10
$arr[0,0] = -16.98665313593701
$arr[0,1] = 2.227915008146631
$arr[0,2] = -16.98665313593701
$arr[0,3] = -16.98665313593701
$arr[0,4] = -16.98665313593701
$arr[0,5] = -16.98665313593701
$arr[0,6] = -16.98665313593701
$arr[0,7] = -16.98665313593701
$arr[0,8] = -16.98665313593701
$arr[0,9] = -16.98665313593701
$arr[1,0] = -16.98665313593701
$arr[1,1] = -16.98665313593701
$arr[1,2] = -16.98665313593701
$arr[1,3] = -16.98665313593701
$arr[1,4] = -16.98665313593701
$arr[1,5] = -16.98665313593701
$arr[1,6] = -16.98665313593701
$arr[1,7] = 2.846534998960875
$arr[1,8] = -16.98665313593701
$arr[1,9] = -16.98665313593701
$arr[2,0] = -16.98665313593701
$arr[2,1] = 2.190924558014966
$arr[2,2] = -16.98665313593701
$arr[2,3] = -16.98665313593701
$arr[2,4] = 3.336904159054071
$arr[2,5] = -16.98665313593701
$arr[2,6] = 1.561163959505361
$arr[2,7] = -16.98665313593701
$arr[2,8] = 2.126491275896428
$arr[2,9] = -16.98665313593701
$arr[3,0] = -16.98665313593701
$arr[3,1] = 2.243669925415227
$arr[3,2] = -16.98665313593701
$arr[3,3] = -16.98665313593701
$arr[3,4] = -16.98665313593701
$arr[3,5] = 2.878368723352323
$arr[3,6] = -16.98665313593701
$arr[3,7] = -16.98665313593701
$arr[3,8] = -16.98665313593701
$arr[3,9] = -16.98665313593701
$arr[4,0] = -16.98665313593701
$arr[4,1] = -16.98665313593701
$arr[4,2] = -16.98665313593701
$arr[4,3] = -16.98665313593701
$arr[4,4] = -16.98665313593701
$arr[4,5] = -16.98665313593701
$arr[4,6] = -16.98665313593701
$arr[4,7] = -16.98665313593701
$arr[4,8] = -16.98665313593701
$arr[4,9] = -16.98665313593701
$arr[5,0] = 2.587572649462829
$arr[5,1] = -16.98665313593701
$arr[5,2] = -16.98665313593701
$arr[5,3] = -16.98665313593701
$arr[5,4] = -16.98665313593701
$arr[5,5] = -16.98665313593701
$arr[5,6] = -16.98665313593701
$arr[5,7] = -16.98665313593701
$arr[5,8] = -16.98665313593701
$arr[5,9] = -16.98665313593701
$arr[6,0] = -16.98665313593701
$arr[6,1] = -16.98665313593701
$arr[6,2] = -16.98665313593701
$arr[6,3] = 1.843629946621229
$arr[6,4] = -16.98665313593701
$arr[6,5] = -16.98665313593701
$arr[6,6] = -16.98665313593701
$arr[6,7] = -16.98665313593701
$arr[6,8] = -16.98665313593701
$arr[6,9] = -16.98665313593701
$arr[7,0] = 3.806220545242185
$arr[7,1] = -16.98665313593701
$arr[7,2] = -16.98665313593701
$arr[7,3] = -16.98665313593701
$arr[7,4] = -16.98665313593701
$arr[7,5] = -16.98665313593701
$arr[7,6] = -16.98665313593701
$arr[7,7] = -16.98665313593701
$arr[7,8] = -16.98665313593701
$arr[7,9] = -16.98665313593701
$arr[8,0] = -16.98665313593701
$arr[8,1] = -16.98665313593701
$arr[8,2] = -16.98665313593701
$arr[8,3] = -16.98665313593701
$arr[8,4] = -16.98665313593701
$arr[8,5] = -16.98665313593701
$arr[8,6] = -16.98665313593701
$arr[8,7] = 1.20381299701101
$arr[8,8] = -16.98665313593701
$arr[8,9] = 1.987618597937753
$arr[9,0] = -16.98665313593701
$arr[9,1] = -16.98665313593701
$arr[9,2] = -16.98665313593701
$arr[9,3] = 3.020298824323739
$arr[9,4] = -16.98665313593701
$arr[9,5] = 2.748375597118064
$arr[9,6] = -16.98665313593701
$arr[9,7] = -16.98665313593701
$arr[9,8] = -16.98665313593701
$arr[9,9] = 1.80853994458124
$arr[10,0] = -16.98665313593701
$arr[10,1] = -16.98665313593701
$arr[10,2] = -16.98665313593701
$arr[10,3] = -16.98665313593701
$arr[10,4] = -16.98665313593701
$arr[10,5] = -16.98665313593701
$arr[10,6] = -16.98665313593701
$arr[10,7] = -16.98665313593701
$arr[10,8] = -16.98665313593701
$arr[10,9] = -16.98665313593701
$arr[11,0] = -16.98665313593701
$arr[11,1] = -16.98665313593701
$arr[11,2] = -16.98665313593701
$arr[11,3] = 2.052097099078769
$arr[11,4] = -16.98665313593701
$arr[11,5] = -16.98665313593701
$arr[11,6] = -16.98665313593701
$arr[11,7] = -16.98665313593701
$arr[11,8] = 2.127457733665634
$arr[11,9] = 1.883076242753027
$arr[12,0] = -16.98665313593701
$arr[12,1] = -16.98665313593701
$arr[12,2] = -16.98665313593701
$arr[12,3] = -16.98665313593701
$arr[12,4] = -16.98665313593701
$arr[12,5] = -16.98665313593701
$arr[12,6] = -16.98665313593701
$arr[12,7] = -16.98665313593701
$arr[12,8] = -16.98665313593701
$arr[12,9] = 2.16618871380603
$arr[13,0] = -16.98665313593701
$arr[13,1] = -16.98665313593701
$arr[13,2] = -16.98665313593701
$arr[13,3] = -16.98665313593701
$arr[13,4] = -16.98665313593701
$arr[13,5] = -16.98665313593701
$arr[13,6] = -16.98665313593701
$arr[13,7] = -16.98665313593701
$arr[13,8] = -16.98665313593701
$arr[13,9] = -16.98665313593701
$arr[14,0] = -16.98665313593701
$arr[14,1] = -16.98665313593701
$arr[14,2] = -16.98665313593701
$arr[14,3] = -16.98665313593701
$arr[14,4] = -16.98665313593701
$arr[14,5] = -16.98665313593701
$arr[14,6] = -16.98665313593701
$arr[14,7] = -16.98665313593701
$arr[14,8] = 2.177185688802528
$arr[14,9] = -16.98665313593701
$arr[15,0] = -16.98665313593701
$arr[15,1] = 1.753657972737565
$arr[15,2] = -16.98665313593701
$arr[15,3] = -16.98665313593701
$arr[15,4] = -16.98665313593701
$arr[15,5] = -16.98665313593701
$arr[15,6] = 1.246765257942578
$arr[15,7] = 1.892921139877956
$arr[15,8] = 1.875829856810221
$arr[15,9] = -16.98665313593701
$arr[16,0] = -16.98665313593701
$arr[16,1] = -16.98665313593701
$arr[16,2] = -16.98665313593701
$arr[16,3] = -16.98665313593701
$arr[16,4] = -16.98665313593701
$arr[16,5] = -16.98665313593701
$arr[16,6] = 1.553917653180798
$arr[16,7] = 1.120528174303961
$arr[16,8] = 1.618728715633626
$arr[16,9] = -16.98665313593701
$arr[17,0] = -16.98665313593701
$arr[17,1] = -16.98665313593701
$arr[17,2] = 4.321917541117664
$arr[17,3] = -16.98665313593701
$arr[17,4] = -16.98665313593701
$arr[17,5] = -16.98665313593701
$arr[17,6] = 1.621836330179098
$arr[17,7] = 1.466975201191672
$arr[17,8] = -16.98665313593701
$arr[17,9] = -16.98665313593701
$arr[18,0] = -16.98665313593701
$arr[18,1] = 0.3973295314862299
$arr[18,2] = -16.98665313593701
$arr[18,3] = -16.98665313593701
$arr[18,4] = 3.306774730250198
$arr[18,5] = -16.98665313593701
$arr[18,6] = 2.047643430479084
$arr[18,7] = 0.8884302357637119
$arr[18,8] = -16.98665313593701
$arr[18,9] = 2.122414863074122
$arr[19,0] = -16.98665313593701
$arr[19,1] = 0.4050687436222234
$arr[19,2] = -16.98665313593701
$arr[19,3] = 2.053221754911876
$arr[19,4] = -16.98665313593701
$arr[19,5] = 2.567322946831185
$arr[19,6] = 2.178805947122352
$arr[19,7] = -16.98665313593701
$arr[19,8] = -16.98665313593701
$arr[19,9] = -16.98665313593701

$ws.Range("B2:K21").Value = $arr
